# Edit: add Week #15 sprint-retrospective bullet points and start Week #16
# (per commit message: "updated week #15 and inputted week #16")
#
# The document has a Week #15 section ending with a bullet "Continuing work
# on Final Project" immediately followed (in the XML) by an otherwise-empty
# bullet paragraph that carries the hidden "_GoBack" bookmark. We insert five
# new bulleted paragraphs between those two paragraphs, being careful not to
# disturb the bookmark paragraph itself.

$d = $word.ActiveDocument

# Locate the anchor paragraph ("Continuing work on Final Project") robustly
# via Find, then translate that hit into a 1-based Paragraphs() index.
$hitRange = $d.Content
[void]$hitRange.Find.Execute("Continuing work on Final Project", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$anchorStart = $hitRange.Paragraphs(1).Range.Start

$anchorIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Start -eq $anchorStart) {
        $anchorIndex = $i
        break
    }
}

# All new paragraphs get inserted right before the paragraph that currently
# follows the anchor (the bookmark-only paragraph). Using InsertParagraphBefore
# on that paragraph repeatedly -- and always re-fetching it by index -- keeps
# pushing it further down while leaving it (and everything before the anchor)
# completely untouched.
$targetIndex = $anchorIndex + 1

# 1) "First Sprint Retrospective on" + " Final Project" (two separate runs,
#    matching the original document's pattern of split runs).
$bp = $d.Paragraphs($targetIndex)
[void]$bp.Range.InsertParagraphBefore()
$newP = $d.Paragraphs($targetIndex)
$xml1 = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="7"/></w:numPr></w:pPr><w:r><w:t>First Sprint Retrospective on</w:t></w:r><w:r><w:t xml:space="preserve"> Final Project</w:t></w:r></w:p>'
[void]$newP.Range.InsertXML($xml1)
$targetIndex = $targetIndex + 1

# 2) "Second big part of Project"
$bp = $d.Paragraphs($targetIndex)
[void]$bp.Range.InsertParagraphBefore()
$newP = $d.Paragraphs($targetIndex)
[void]$newP.Range.InsertAfter("Second big part of Project")
$targetIndex = $targetIndex + 1

# 3) "Joined a new Group"
$bp = $d.Paragraphs($targetIndex)
[void]$bp.Range.InsertParagraphBefore()
$newP = $d.Paragraphs($targetIndex)
[void]$newP.Range.InsertAfter("Joined a new Group")
$targetIndex = $targetIndex + 1

# 4) "Working on new content for project"
$bp = $d.Paragraphs($targetIndex)
[void]$bp.Range.InsertParagraphBefore()
$newP = $d.Paragraphs($targetIndex)
[void]$newP.Range.InsertAfter("Working on new content for project")
$targetIndex = $targetIndex + 1

# 5) "Week #16" -- top-level bullet (ilvl 0) rather than the nested sub-bullet
#    level (ilvl 1) used by the other new paragraphs.
$bp = $d.Paragraphs($targetIndex)
[void]$bp.Range.InsertParagraphBefore()
$newP = $d.Paragraphs($targetIndex)
[void]$newP.Range.InsertAfter("Week #16")
$newP.Range.ListFormat.ListLevelNumber = 1
$targetIndex = $targetIndex + 1
